$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 2945
$ws.Cells.Item(3, 6).Value = 6409
$ws.Cells.Item(4, 6).Value = 2517
$ws.Cells.Item(6, 6).Value = 525
$ws.Cells.Item(7, 6).Value = 54
$ws.Cells.Item(9, 6).Value = 2925
$ws.Cells.Item(10, 6).Value = 351
$ws.Cells.Item(12, 6).Value = 7347
$ws.Cells.Item(13, 6).Value = 333
$ws.Cells.Item(16, 6).Value = 243
$ws.Cells.Item(19, 6).Value = 8946
$ws.Cells.Item(27, 6).Value = 104
$ws.Cells.Item(28, 6).Value = 35
$ws.Cells.Item(30, 6).Value = 45
$ws.Cells.Item(31, 6).Value = 65
$ws.Cells.Item(36, 6).Value = 40
$ws.Cells.Item(38, 6).Value = 743
$ws.Cells.Item(39, 6).Value = 3863
$ws.Cells.Item(40, 6).Value = 200
$ws.Cells.Item(41, 6).Value = 30
$ws.Cells.Item(43, 6).Value = 54
$ws.Cells.Item(44, 6).Value = 16
$ws.Cells.Item(49, 6).Value = 46

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 47
$ws.Cells.Item(12, 6).Value = 2
$ws.Cells.Item(16, 6).Value = 167
$ws.Cells.Item(18, 6).Value = 33

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 47
$ws.Cells.Item(3, 6).Value = 2945
$ws.Cells.Item(6, 6).Value = 6409
$ws.Cells.Item(7, 6).Value = 2517
$ws.Cells.Item(10, 6).Value = 525
$ws.Cells.Item(11, 6).Value = 54
$ws.Cells.Item(13, 6).Value = 2925
$ws.Cells.Item(14, 6).Value = 351
$ws.Cells.Item(18, 6).Value = 7347
$ws.Cells.Item(19, 6).Value = 333
$ws.Cells.Item(22, 6).Value = 243
$ws.Cells.Item(24, 6).Value = 8946
$ws.Cells.Item(30, 6).Value = 104
$ws.Cells.Item(31, 6).Value = 45
$ws.Cells.Item(32, 6).Value = 65
$ws.Cells.Item(37, 6).Value = 40
$ws.Cells.Item(39, 6).Value = 743
$ws.Cells.Item(40, 6).Value = 167
$ws.Cells.Item(41, 6).Value = 3863
$ws.Cells.Item(42, 6).Value = 200
$ws.Cells.Item(43, 6).Value = 30
$ws.Cells.Item(49, 6).Value = 46

